# "ready for iteration 3" -- add a new gemma/qwen/llama model comparison row,
# drop the now-unused raw per-rater columns C and L on the lower table, and
# refresh a couple of literal scores plus some cosmetic view/column-width state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------------
# 1. New row 12: extra model-name header cells above the per-exam detail rows
# ---------------------------------------------------------------------------
$ws.Range("D12").Value = "gemma3 4b"
$ws.Range("E12").Value = "qwen 3B"
$ws.Range("G12").Value = "gemma"
$ws.Range("H12").Value = "qwen"
$ws.Range("K12").Value = "gemma"
$ws.Range("L12").Value = "llama"
$ws.Range("M12").Value = "qwen"

# ---------------------------------------------------------------------------
# 2. Drop the old per-rater raw columns C and L for rows 13:19 (the running
#    AVERAGE formulas in B/F/J pick up the change automatically).
# ---------------------------------------------------------------------------
$ws.Range("C13:C19").ClearContents()
$ws.Range("L13:L19").ClearContents()

# ---------------------------------------------------------------------------
# 3. A couple of literal score corrections on row 15.
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = 0.95
$ws.Range("M15").Value = 1

# ---------------------------------------------------------------------------
# 4. Column width / visibility touch-ups (cosmetic formatting state).
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 5.333333333333333
$ws.Columns("D").ColumnWidth = 8
$ws.Columns("E").ColumnWidth = 10.666666666666666

$ws.Columns("G").ColumnWidth = 5.5
$ws.Columns("H").ColumnWidth = 9
$ws.Columns("I").ColumnWidth = 5.166666666666667

$ws.Columns("L").ColumnWidth = 5
$ws.Columns("M").ColumnWidth = 7
$ws.Columns("N").ColumnWidth = 5.166666666666667

$ws.Columns("P").Hidden = $false
$ws.Columns("P").ColumnWidth = 4.666666666666667
$ws.Columns("Q").Hidden = $false
$ws.Columns("Q").ColumnWidth = 7.666666666666667
$ws.Columns("R").Hidden = $false
$ws.Columns("R").ColumnWidth = 10.333333333333334

# ---------------------------------------------------------------------------
# 5. Refresh the active selection to match where editing left off.
# ---------------------------------------------------------------------------
$ws.Range("P17").Select()
